# Apply cryptocurrency price/volume updates to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D2").Value = "68.395.32"
$ws.Range("E2").Value = "  +1.96%  "
$ws.Range("D3").Value = "2.641.29"
$ws.Range("E3").Value = "  +1.30%  "
$ws.Range("E4").Value = "  -0.16%  "
$ws.Range("D5").Value = "'599.72"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.44%  "
$ws.Range("D6").Value = "'154.15"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.74%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").Value = "'0.545"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.18%  "
$ws.Range("D9").Value = "2.640.16"
$ws.Range("E9").Value = "  +1.27%  "
$ws.Range("E10").Value = "  +11.97%  "
$ws.Range("E11").Value = "  -0.46%  "
$ws.Range("E12").Value = "  +0.73%  "
$ws.Range("E13").Value = "  +0.41%  "
$ws.Range("D14").Value = "'27.64"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.45%  "
$ws.Range("E15").Value = "  +5.52%  "
$ws.Range("D16").Value = "3.125.20"
$ws.Range("E16").Value = "  +1.43%  "
$ws.Range("D17").Value = "68.286.88"
$ws.Range("E17").Value = "  +2.07%  "
$ws.Range("D18").Value = "2.629.16"
$ws.Range("E18").Value = "  +0.79%  "
$ws.Range("D19").Value = "'11.37"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.33%  "
$ws.Range("D20").Value = "'366.76"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.12%  "
$ws.Range("D21").Value = "'7.44"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.19%  "
$ws.Range("D22").Value = "'4.24"
$ws.Range("D22").Style = "Normal"
$ws.Range("E23").Value = "  +0.12%  "
$ws.Range("D24").Value = "'2.09"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.84%  "
$ws.Range("D25").Value = "'72.90"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +7.88%  "
$ws.Range("E26").Value = "  +0.04%  "
$ws.Range("E27").Value = "  -1.94%  "
$ws.Range("D28").Value = "2.764.19"
$ws.Range("E28").Value = "  +0.87%  "
$ws.Range("E29").Value = "  +3.95%  "
$ws.Range("E30").Value = "  +0.03%  "
$ws.Range("D31").Value = "'573.70"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.39%  "
$ws.Range("D32").Value = "'7.93"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.38%  "
$ws.Range("E33").Value = "  +2.01%  "
$ws.Range("E34").Value = "  +2.86%  "
$ws.Range("D35").Value = "'1.00"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.09%  "
$ws.Range("E36").Value = "  +4.41%  "
$ws.Range("E37").Value = "  +1.28%  "
$ws.Range("D38").Value = "'159.91"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.29%  "
$ws.Range("D39").Value = "'19.18"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.76%  "
$ws.Range("D40").Value = "'1.89"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.48%  "
$ws.Range("E41").Value = "  +0.26%  "
$ws.Range("D42").Value = "'5.34"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.20%  "
$ws.Range("D43").Value = "'2.64"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.18%  "
$ws.Range("D44").Value = "'17.62"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +5.10%  "
$ws.Range("E45").Value = "  +10.80%  "
$ws.Range("E46").Value = "  +0.08%  "
$ws.Range("D47").Value = "'40.47"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.91%  "
$ws.Range("D48").Value = "'155.61"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.65%  "
$ws.Range("E49").Value = "  -0.13%  "
$ws.Range("D50").Value = "'21.93"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.07%  "
$ws.Range("E51").Value = "  +0.38%  "
